$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '48.246.87'
$ws.Range('E2').Value = '  +2.15%  '
$ws.Range('D3').Value = '2.525.48'
$ws.Range('E3').Value = '  +1.38%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '323.76'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.49%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '109.66'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.81%  '
$ws.Range('E7').Value = '  +1.09%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.556'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +3.79%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.02'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +5.45%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.49'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +11.72%  '
$ws.Range('E12').Value = '  +1.69%  '
$ws.Range('E13').Value = '  +1.25%  '
$ws.Range('E14').Value = '  +1.78%  '
$ws.Range('D15').Value = '2.922.16'
$ws.Range('E15').Value = '  +1.53%  '
$ws.Range('D16').Value = '2.522.91'
$ws.Range('E16').Value = '  +1.20%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.859'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +1.33%  '
$ws.Range('D18').Value = '48.074.81'
$ws.Range('E18').Value = '  +2.00%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.27'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +3.89%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.64'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.24%  '
$ws.Range('D21').Value = '0.0₃0950'
$ws.Range('E21').Value = '  +1.31%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.72'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '72.18'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +2.15%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '268.29'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +8.65%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.58'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.11%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.25'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.52%  '
$ws.Range('E27').Value = '  -0.29%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.15'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.92%  '
$ws.Range('E29').Value = '  +2.90%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.22'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.62%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '36.00'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.74%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '49.68'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.51%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.04'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.29%  '
$ws.Range('E34').Value = '  -0.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0795'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.11%  '
$ws.Range('E37').Value = '  +1.61%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.76'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.32%  '
$ws.Range('E39').Value = '  +1.58%  '
$ws.Range('E40').Value = '  +0.40%  '
$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '120.41'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.25%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '22.05'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +4.18%  '
$ws.Range('E43').Value = '  -1.66%  '
$ws.Range('E44').Value = '  +1.92%  '
$ws.Range('D45').Value = '2.024.09'
$ws.Range('E45').Value = '  +1.56%  '
$ws.Range('E46').Value = '  +4.80%  '
$ws.Range('E47').Value = '  +6.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.03'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.52%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.18'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.81%  '
$ws.Range('E50').Value = '  +3.06%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '79.76'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +2.80%  '
